# Update latest output (run 55)
$wb = $excel.ActiveWorkbook

# --- Schedule sheet ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 878.5684342500001
$schedule.Range("F2").Value = 14.52659448164683

# --- Detailed sheet ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B14").Value = 56.98
$detailed.Range("B15").Value = 50.46272

$detailed.Range("C17").Value = "historical"

$detailed.Range("B18").Value = 12.94586
$detailed.Range("C18").Value = "historical"

$detailed.Range("B19").Value = 18.67563
$detailed.Range("B20").Value = 27.87441
$detailed.Range("B22").Value = 33.86823
$detailed.Range("B23").Value = 35.88
$detailed.Range("B24").Value = 36.06092
$detailed.Range("B25").Value = 44.06468
$detailed.Range("B26").Value = 36.06046
$detailed.Range("B27").Value = 43.1301
$detailed.Range("B28").Value = 46.34072
$detailed.Range("B30").Value = 34.78031
$detailed.Range("B32").Value = 36.06033
$detailed.Range("B34").Value = 36.05879
$detailed.Range("B35").Value = 16.07216
$detailed.Range("B36").Value = -0.41405
$detailed.Range("B37").Value = -2.99806
$detailed.Range("B38").Value = -2.91785
$detailed.Range("B39").Value = -2.88418
$detailed.Range("B40").Value = 7.52239
$detailed.Range("B41").Value = 32.87797
$detailed.Range("B42").Value = 56.40935
$detailed.Range("B43").Value = 10.22525
$detailed.Range("B44").Value = 19.64731
$detailed.Range("B45").Value = 22.87055
$detailed.Range("B46").Value = 36.06045
$detailed.Range("B49").Value = 56.98
